# Applies factsheet text-formatting edits described in the commit:
# "Update factsheets with text edits from COMM"
# Converts numeric filer counts to comma-formatted text, fixes the
# Alpine County placeholder row, and appends a County "Total" row.
$wb = $excel.ActiveWorkbook

# --- Overall: No. of 990 Filers total, numeric -> comma text ---
$ws = $wb.Worksheets.Item("Overall")
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "10,513"

# --- County: column B numeric -> comma text ---
$ws = $wb.Worksheets.Item("County")
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "873"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "13"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "71"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "13"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "3"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "286"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "10"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "47"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "168"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "4"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "96"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "21"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "12"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "106"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "13"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "16"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "7"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "2,519"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "14"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "241"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "11"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "68"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "24"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "8"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "11"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "122"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "76"
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "73"
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "587"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "113"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "14"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "281"
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "463"
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "17"
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "266"
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "818"
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "858"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "96"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "113"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "255"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "209"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "502"
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "151"
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "65"
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "5"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "16"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "59"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "228"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "65"
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "17"
$ws.Range("B52").NumberFormat = "@"
$ws.Range("B52").Value = "11"
$ws.Range("B53").NumberFormat = "@"
$ws.Range("B53").Value = "9"
$ws.Range("B54").NumberFormat = "@"
$ws.Range("B54").Value = "63"
$ws.Range("B55").NumberFormat = "@"
$ws.Range("B55").Value = "20"
$ws.Range("B56").NumberFormat = "@"
$ws.Range("B56").Value = "210"
$ws.Range("B57").NumberFormat = "@"
$ws.Range("B57").Value = "67"
$ws.Range("B58").NumberFormat = "@"
$ws.Range("B58").Value = "9"

# --- Congressional District: column B numeric -> comma text ---
$ws = $wb.Worksheets.Item("Congressional District")
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "210"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "176"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "823"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "729"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "56"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "126"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "238"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "265"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "158"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "213"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "258"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "534"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "115"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "143"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "55"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "100"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "351"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "73"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "175"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "84"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "231"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "108"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "256"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "294"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "94"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "207"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "90"
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "383"
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "99"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "306"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "279"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "56"
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "84"
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "285"
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "116"
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "112"
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "235"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "83"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "104"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "85"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "156"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "197"
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "112"
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "149"
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "178"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "311"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "217"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "104"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "198"
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "256"
$ws.Range("B52").NumberFormat = "@"
$ws.Range("B52").Value = "152"
$ws.Range("B53").NumberFormat = "@"
$ws.Range("B53").Value = "94"
$ws.Range("B54").NumberFormat = "@"
$ws.Range("B54").Value = "10,513"

# --- Size: column B numeric -> comma text ---
$ws = $wb.Worksheets.Item("Size")
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2,809"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "3,237"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "1,720"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "867"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "1,511"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "369"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "10,513"

# --- Subsector: column B numeric -> comma text ---
$ws = $wb.Worksheets.Item("Subsector")
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1,192"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1,432"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "559"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "867"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "54"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "3,022"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "159"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "2"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "860"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "231"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "2,031"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "104"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "10,513"

# --- County: row 59 (Alpine County) placeholder values reformatted ---
$ws = $wb.Worksheets.Item("County")
$ws.Range("B59").NumberFormat = "@"
$ws.Range("B59").Value = "0.00%"
$ws.Range("C59").NumberFormat = "@"
$ws.Range("C59").Value = "$0"
$ws.Range("D59").NumberFormat = "@"
$ws.Range("D59").Value = "0.00%"
$ws.Range("E59").NumberFormat = "@"
$ws.Range("E59").Value = "0.00%"
$ws.Range("F59").NumberFormat = "@"
$ws.Range("F59").Value = "0.00%"

# --- County: append new Total row (row 60) ---
$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value = "Total"
$ws.Range("B60").NumberFormat = "@"
$ws.Range("B60").Value = "10,513"
$ws.Range("C60").NumberFormat = "@"
$ws.Range("C60").Value = "$42,132,825,415"
$ws.Range("D60").NumberFormat = "@"
$ws.Range("D60").Value = "9.03%"
$ws.Range("E60").NumberFormat = "@"
$ws.Range("E60").Value = "-13.04%"
$ws.Range("F60").NumberFormat = "@"
$ws.Range("F60").Value = "65.68%"

